# Apply the "Alvearie FHIR IG" gh-pages deployment update to the
# StructureDefinition-cost-sharing-reduction-variant workbook.
#
# Sheet 1 ("Metadata"): bump Version/Date, set Publisher, replace the
# duplicated "Contact" row with a new "Jurisdiction" row (net: 21 -> 20 rows).
# Sheet 2 ("Elements"): update the Short/Definition text of the root
# Extension element row to reflect the profile's own title/description.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# Remove the duplicate "Contact" / "No display for ContactDetail" row
# (originally row 11) -- everything below shifts up by one row.
$ws1.Rows.Item(11).Delete()

# Update metadata values on the "Metadata" sheet.
$ws1.Range("B3").Value = "6.0.0"
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$ws1.Range("B9").Value = "Alvearie Team"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

$ws2 = $wb.Worksheets.Item(2)

# Update the Short / Definition columns for the root Extension row.
$ws2.Range("K2").Value = "Cost Sharing Reduction Variant"
$ws2.Range("L2").Value = "Code for Affordable Care Act (ACA) cost sharing reduction variant of the associated plan"
